$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.04780798405408859
$ws.Range("B2").Value = 2.907422416370558
$ws.Range("C2").Value = 3.823013305664062
$ws.Range("D2").Value = 8
